# Applies the 2019-09-11-oop.pptx edit:
#  1. Update every cached "datetimeFigureOut" field (slide master, all
#     slide layouts, notes master) from "9/11/2019" to "9/12/19".
#  2. Re-word the agenda line on slide 1 from
#       "Agenda for Monday, September 11"
#     to three runs reading
#       "Agenda " / "for Wednesday, " / "September 11"
#     (the following "h" superscript + " from 2 to 2:50pm CST:" runs are
#     left untouched).

$p = $ppt.ActivePresentation

function Set-DatePlaceholderText {
    param($shapes, [string]$newText)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $isDate = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq 16) {
                $isDate = $true
            }
        } catch {
            $isDate = $false
        }
        if ($isDate) {
            $sh.TextFrame.TextRange.Text = $newText
        }
    }
}

# --- 1. Date fields -------------------------------------------------

# Slide master.
Set-DatePlaceholderText -shapes $p.SlideMaster.Shapes -newText "9/12/19"

# Every slide layout (custom layout) hanging off the (single) design.
$design = $p.Designs.Item(1)
$layouts = $design.SlideMaster.CustomLayouts
for ($j = 1; $j -le $layouts.Count; $j++) {
    $cl = $layouts.Item($j)
    Set-DatePlaceholderText -shapes $cl.Shapes -newText "9/12/19"
}

# Notes master.
Set-DatePlaceholderText -shapes $p.NotesMaster.Shapes -newText "9/12/19"

# --- 2. Slide 1 agenda line ------------------------------------------

$slide1 = $p.Slides.Item(1)
$contentShape = $slide1.Shapes.Item(2)
$tr = $contentShape.TextFrame.TextRange
$para1 = $tr.Paragraphs(1, 1)

# Original paragraph 1 text: "Agenda for Monday, September 11h from 2 to 2:50pm CST:"
# Characters 1-31 hold "Agenda for Monday, September 11"; split that run
# into three pieces by rewriting the middle slice in place - this keeps
# the "h" (superscript) and trailing " from 2 to 2:50pm CST:" runs intact.
$middle = $para1.Characters(8, 12)   # "for Monday, "
$middle.Text = "for Wednesday, "

Write-Output "Done."
